# Estadisticos Matutinos 15 Oct
# Fills in the "Reprobados/Aprobados/Por_Apro/Promedio" stats for the two
# partial-period sheets + the final sheet, and adds the list of rescatable
# (make-up exam) students to the "Rescatables" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Estadisticos 1P" and "Estadisticos Final": same figures in both sheets.
#    Columns: A=Mat B=Grupo C=Totales D=Blancos E=Reprobados F=Aprobados
#             G=Por_Apro H=Promedio
# ---------------------------------------------------------------------------
$statsRows = @(
    @{ Row = 2; Blancos = 9;  Aprobados = 25; PorApro = 73.53; Promedio = 7.3 },
    @{ Row = 3; Blancos = 9;  Aprobados = 22; PorApro = 70.97; Promedio = 7.4 },
    @{ Row = 4; Blancos = 7;  Aprobados = 29; PorApro = 80.56; Promedio = 7.6 },
    @{ Row = 5; Blancos = 11; Aprobados = 32; PorApro = 74.42; Promedio = 7.5 },
    @{ Row = 6; Blancos = 11; Aprobados = 33; PorApro = 75;    Promedio = 7.4 },
    @{ Row = 7; Blancos = 1;  Aprobados = 23; PorApro = 95.83; Promedio = 7.3 }
)

foreach ($sheetName in @("Estadisticos 1P", "Estadisticos Final")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($r in $statsRows) {
        $row = $r.Row
        $ws.Cells.Item($row, 4).Value = $r.Blancos
        $ws.Cells.Item($row, 6).Value = $r.Aprobados
        $ws.Cells.Item($row, 7).Value = $r.PorApro
        $ws.Cells.Item($row, 8).Value = $r.Promedio
    }
}

# ---------------------------------------------------------------------------
# 2) "Estadisticos 2P": only the Reprobados (E) column gets filled in.
# ---------------------------------------------------------------------------
$reprobadosRows = @(
    @{ Row = 2; Reprobados = 25 },
    @{ Row = 3; Reprobados = 22 },
    @{ Row = 4; Reprobados = 29 },
    @{ Row = 5; Reprobados = 32 },
    @{ Row = 6; Reprobados = 33 },
    @{ Row = 7; Reprobados = 23 }
)

$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
foreach ($r in $reprobadosRows) {
    $ws2.Cells.Item($r.Row, 5).Value = $r.Reprobados
}

# ---------------------------------------------------------------------------
# 3) "Rescatables": add the list of students who need a make-up exam.
#    Columns: A=NC B=Paterno C=Materno D=Nombres E=Nombre_Largo F=Grupo G=Reprobadas
# ---------------------------------------------------------------------------
$rescatables = @(
    @{ Row = 2; NC = 21330051920241; Paterno = "MONTERO";   Materno = "LOPEZ";     Nombres = "MARIA DEL PILAR"; Grupo = "1CM" },
    @{ Row = 3; NC = 21330051920297; Paterno = "DOMINGUEZ"; Materno = "APALE";     Nombres = "EDNA";            Grupo = "1EM" },
    @{ Row = 4; NC = 21330051920323; Paterno = "ROMERO";    Materno = "SANCHEZ";   Nombres = "DIEGO EMILIO";    Grupo = "1EM" },
    @{ Row = 5; NC = 21330051920303; Paterno = "HUERTA";    Materno = "GONZALEZ";  Nombres = "YULIET";          Grupo = "1EM" },
    @{ Row = 6; NC = 21330051920308; Paterno = "MOTA";      Materno = "CERON";     Nombres = "ANGEL DIEGO";     Grupo = "1EM" }
)

$nombreLargo = "TECNOLOGÍAS DE LA INFORMACIÓN Y LA COMUNICACIÓN"

$ws4 = $wb.Worksheets.Item("Rescatables")

# Shared strings are interned in first-seen order, so write column-by-column
# (Paterno x5, then Materno x5, then Nombres x5, ...) to reproduce the same
# sharedStrings.xml ordering as the source workbook.
foreach ($r in $rescatables) {
    $ws4.Cells.Item($r.Row, 1).Value = $r.NC
}
foreach ($r in $rescatables) {
    $ws4.Cells.Item($r.Row, 2).Value = $r.Paterno
}
foreach ($r in $rescatables) {
    $ws4.Cells.Item($r.Row, 3).Value = $r.Materno
}
foreach ($r in $rescatables) {
    $ws4.Cells.Item($r.Row, 4).Value = $r.Nombres
}
foreach ($r in $rescatables) {
    $ws4.Cells.Item($r.Row, 5).Value = $nombreLargo
}
foreach ($r in $rescatables) {
    $ws4.Cells.Item($r.Row, 6).Value = $r.Grupo
}
foreach ($r in $rescatables) {
    $ws4.Cells.Item($r.Row, 7).Value = 6
}
